$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.184.01"
$ws.Range("E2").Value = "  +5.02%  "
$ws.Range("D3").Value = "4.076.08"
$ws.Range("E3").Value = "  +5.38%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "523.43"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.21"
$ws.Range("E6").Value = "  +3.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.719"
$ws.Range("E7").Value = "  +18.70%  "
$ws.Range("D8").Value = "4.067.90"
$ws.Range("E8").Value = "  +5.39%  "
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.776"
$ws.Range("E10").Value = "  +9.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.181"
$ws.Range("E11").Value = "  +7.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000338"
$ws.Range("E12").Value = "  +4.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "48.78"
$ws.Range("E13").Value = "  +16.67%  "
$ws.Range("E14").Value = "  +8.65%  "
$ws.Range("D15").Value = "4.727.81"
$ws.Range("E15").Value = "  +5.63%  "
$ws.Range("D16").Value = "4.071.96"
$ws.Range("E16").Value = "  +6.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.53"
$ws.Range("E17").Value = "  +2.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "21.34"
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.24"
$ws.Range("D21").Value = "72.275.91"
$ws.Range("E21").Value = "  +5.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "449.90"
$ws.Range("E22").Value = "  +6.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "103.62"
$ws.Range("E23").Value = "  +18.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.64"
$ws.Range("E24").Value = "  +7.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "15.03"
$ws.Range("E25").Value = "  +7.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.09"
$ws.Range("E26").Value = "  +2.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.52"
$ws.Range("E27").Value = "  +1.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.20"
$ws.Range("E28").Value = "  +5.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.09"
$ws.Range("E29").Value = "  +5.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.82"
$ws.Range("E30").Value = "  +2.69%  "
$ws.Range("E31").Value = "  +16.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.71"
$ws.Range("E32").Value = "  +4.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.133"
$ws.Range("E33").Value = "  +5.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "684.57"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "67.96"
$ws.Range("E35").Value = "  +1.20%  "
$ws.Range("E36").Value = "  +13.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.60"
$ws.Range("E37").Value = "  +6.88%  "
$ws.Range("D38").Value = "0.0₃0897"
$ws.Range("E38").Value = "  +6.05%  "
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.155"
$ws.Range("E40").Value = "  +5.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.50"
$ws.Range("E41").Value = "  +8.94%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("E43").Value = "  +5.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.22"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.157"
$ws.Range("E46").Value = "  +12.44%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.71"
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.93"
$ws.Range("E48").Value = "  +16.80%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000285"
$ws.Range("E50").Value = "  +7.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.08"
$ws.Range("E51").Value = "  +4.99%  "

Write-Output "Applied cryptos list update"
